$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row (row 1) ---------------------------------------------------
# I1: "source" -> "shade_tol"
$ws.Range("I1").Value2 = "shade_tol"

# J1: new header, same look as existing headers (bold, general format)
$ws.Range("J1").Value2 = "inflorescence size"
$ws.Range("J1").Font.Bold = $true

# K1: new header, bold + text ("@") number format
$ws.Range("K1").Value2 = "flower_size"
$ws.Range("K1").Font.Bold = $true
$ws.Range("K1").NumberFormat = "@"

# --- Row 2 ------------------------------------------------------------------
$ws.Range("I2").Value2 = "intolerant"
$ws.Range("K2").Value2 = 5
$ws.Range("K2").NumberFormat = "@"

# L1: new header (bold + text format) then L2 data
$ws.Range("L1").Value2 = "fruit_size"
$ws.Range("L1").Font.Bold = $true
$ws.Range("L1").NumberFormat = "@"

$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value2 = "5-8"

# --- Row 3 --------------------------------------------------------------------
$ws.Range("L3").NumberFormat = "@"
$ws.Range("L3").Value2 = "5-12"

$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value2 = "3-4"

# --- Row 5 (intolerant repeats; written ahead of row 3/4's "tolerant") --------
$ws.Range("I5").Value2 = "intolerant"

# --- Row 4 (fruit size) -------------------------------------------------------
$ws.Range("L4").NumberFormat = "@"
$ws.Range("L4").Value2 = "1"

# --- Row 5 (fruit size) -------------------------------------------------------
$ws.Range("L5").NumberFormat = "@"
$ws.Range("L5").Value2 = "2-3"

# --- Row 3 (shade tolerance, first "tolerant") --------------------------------
$ws.Range("I3").Value2 = "tolerant"

# --- Row 4 (shade tolerance repeat + inflorescence size numeric) -------------
$ws.Range("I4").Value2 = "tolerant"
$ws.Range("J4").Value2 = 5

# --- Row 6 ---------------------------------------------------------------------
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value2 = "1.5-2"

$ws.Range("L6").NumberFormat = "@"
$ws.Range("L6").Value2 = "1-1.5"

$ws.Range("I6").Value2 = "tolerant"

# --- Row 7 -----------------------------------------------------------------------
$ws.Range("I7").Value2 = "intolerant"

$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value2 = "4"

# --- Row 8 -------------------------------------------------------------------------
$ws.Range("J8").Value2 = "8-10,15-20"

# --- View / selection state ----------------------------------------------------------
$ws.Activate()
$ws.Range("K8").Select() | Out-Null
